$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 112241886
$ws.Range("B2").Value = 78713
$ws.Range("Q2").Value = 553306
$ws.Range("R2").Value = 7007600

# Row 3
$ws.Range("A3").Value = 112241888
$ws.Range("B3").Value = 78713
$ws.Range("Q3").Value = 553253
$ws.Range("R3").Value = 7007768

# Row 4
$ws.Range("A4").Value = 112241889
$ws.Range("B4").Value = 78713
$ws.Range("R4").Value = 7007769

# Row 5
$ws.Range("A5").Value = 112241885
$ws.Range("B5").Value = 78713
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 6458
$ws.Range("F5").Value = "Lunglav"
$ws.Range("G5").Value = "Lobaria pulmonaria"
$ws.Range("H5").Value = "(L.) Hoffm."
$ws.Range("Q5").Value = 553321
$ws.Range("R5").Value = 7007611

# Row 6
$ws.Range("A6").Value = 112241887
$ws.Range("B6").Value = 78713
$ws.Range("Q6").Value = 553188
$ws.Range("R6").Value = 7007666

# Row 7
$ws.Range("B7").Value = 78713

# Row 8
$ws.Range("A8").Value = 112241875
$ws.Range("B8").Value = 78740
$ws.Range("D8").Value = "LC"
$ws.Range("E8").Value = 6462
$ws.Range("F8").Value = "Stuplav"
$ws.Range("G8").Value = "Nephroma bellum"
$ws.Range("H8").Value = "(Spreng.) Tuck."
$ws.Range("Q8").Value = 553188
$ws.Range("R8").Value = 7007668
